$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2975.225
$ws.Range("I138").Value = 1813.5714
$ws.Range("J138").Value = 3600.7307
$ws.Range("K138").Value = 5440.7142
$ws.Range("L138").Value = 10802.1921
$ws.Range("M138").Value = -300.7142000000003
$ws.Range("N138").Value = -21082.1921

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H64").Value = 31418.2
$ws.Range("I64").Value = 32000
$ws.Range("J64").Value = 31272.75
$ws.Range("K64").Value = 32000
$ws.Range("L64").Value = 31272.75
$ws.Range("M64").Value = -31752
$ws.Range("N64").Value = -31768.75

$ws.Range("H67").Value = 31418.2
$ws.Range("I67").Value = 32000
$ws.Range("J67").Value = 31272.75
$ws.Range("K67").Value = 32000
$ws.Range("L67").Value = 31272.75
$ws.Range("M67").Value = -31142
$ws.Range("N67").Value = -32988.75

$ws.Range("H74").Value = 8168.7646
$ws.Range("I74").Value = 2523.6667
$ws.Range("J74").Value = 50507
$ws.Range("K74").Value = 2523.6667
$ws.Range("L74").Value = 50507
$ws.Range("M74").Value = -1649.6667
$ws.Range("N74").Value = -52255

$ws.Range("H77").Value = 8168.7646
$ws.Range("I77").Value = 2523.6667
$ws.Range("J77").Value = 50507
$ws.Range("K77").Value = 12618.3335
$ws.Range("L77").Value = 252535
$ws.Range("M77").Value = -8250.333500000001
$ws.Range("N77").Value = -261271

$ws.Range("H92").Value = 48000
$ws.Range("J92").Value = 48000
$ws.Range("L92").Value = 48000
$ws.Range("N92").Value = -52992

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 38181
$ws.Range("J62").Value = 38181
$ws.Range("L62").Value = 38181
$ws.Range("N62").Value = -39553

$ws.Range("H65").Value = 38181
$ws.Range("J65").Value = 38181
$ws.Range("L65").Value = 114543
$ws.Range("N65").Value = -121407

$ws.Range("H92").Value = 25401
$ws.Range("J92").Value = 25401
$ws.Range("L92").Value = 25401
$ws.Range("N92").Value = -30393

$ws.Range("H105").Value = 1362989.6
$ws.Range("I105").Value = 1740258.9
$ws.Range("J105").Value = 4820
$ws.Range("K105").Value = 1740258.9
$ws.Range("L105").Value = 4820
$ws.Range("M105").Value = -1738511.9
$ws.Range("N105").Value = -8314

$ws.Range("H134").Value = 42968.582
$ws.Range("I134").Value = 1357.9131
$ws.Range("K134").Value = 4073.7393
$ws.Range("M134").Value = -1538.7393

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 37867
$ws.Range("J92").Value = 37867
$ws.Range("L92").Value = 37867
$ws.Range("N92").Value = -42859

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2869.8704
$ws.Range("J68").Value = 5491.4346
$ws.Range("L68").Value = 16474.3038
$ws.Range("N68").Value = -18096.3038

$ws.Range("H71").Value = 2869.8704
$ws.Range("J71").Value = 5491.4346
$ws.Range("L71").Value = 49422.9114
$ws.Range("N71").Value = -57534.9114

$ws.Range("H109").Value = 3806
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 3806
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 11418
$ws.Range("M109").ClearContents()
$ws.Range("N109").Value = -13498

$ws.Range("H119").Value = 5507.25
$ws.Range("I119").Value = 1029
$ws.Range("J119").Value = 7000
$ws.Range("K119").Value = 3087
$ws.Range("L119").Value = 21000
$ws.Range("M119").Value = 1751
$ws.Range("N119").Value = -30676

$ws.Range("H131").Value = 1261.1666
$ws.Range("I131").Value = 508.8889
$ws.Range("J131").Value = 1466.3334
$ws.Range("K131").Value = 1526.6667
$ws.Range("L131").Value = 4399.0002
$ws.Range("M131").Value = 3513.3333
$ws.Range("N131").Value = -14479.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 10363.875
$ws.Range("I80").Value = 17835
$ws.Range("J80").Value = 5881.2
$ws.Range("K80").Value = 17835
$ws.Range("L80").Value = 5881.2
$ws.Range("M80").Value = -16837
$ws.Range("N80").Value = -7877.2

$ws.Range("H83").Value = 10363.875
$ws.Range("I83").Value = 17835
$ws.Range("J83").Value = 5881.2
$ws.Range("K83").Value = 89175
$ws.Range("L83").Value = 29406
$ws.Range("M83").Value = -84183
$ws.Range("N83").Value = -39390

$ws.Range("H126").Value = 2947.3684
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 3636.3635
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 10909.0905
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -15849.0905

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2943.4194
$ws.Range("I40").Value = 2809.85
$ws.Range("J40").Value = 3186.2727
$ws.Range("K40").Value = 2809.85
$ws.Range("L40").Value = 3186.2727
$ws.Range("M40").Value = -2673.85
$ws.Range("N40").Value = -3458.2727

$ws.Range("H46").Value = 774.2857
$ws.Range("I46").Value = 671.1111
$ws.Range("J46").Value = 960
$ws.Range("K46").Value = 671.1111
$ws.Range("L46").Value = 960
$ws.Range("M46").Value = -483.1111
$ws.Range("N46").Value = -1336

$ws.Range("H61").Value = 16357.192
$ws.Range("I61").Value = 15881.182
$ws.Range("J61").Value = 18975.25
$ws.Range("K61").Value = 15881.182
$ws.Range("L61").Value = 18975.25
$ws.Range("M61").Value = -15679.182
$ws.Range("N61").Value = -19379.25

$ws.Range("H113").Value = 16357.192
$ws.Range("I113").Value = 15881.182
$ws.Range("J113").Value = 18975.25
$ws.Range("K113").Value = 15881.182
$ws.Range("L113").Value = 18975.25
$ws.Range("M113").Value = -13711.182
$ws.Range("N113").Value = -23315.25

$ws.Range("H122").Value = 5938.788
$ws.Range("I122").Value = 5783.3335
$ws.Range("J122").Value = 7493.3335
$ws.Range("K122").Value = 17350.0005
$ws.Range("L122").Value = 22480.0005
$ws.Range("M122").Value = -14900.0005
$ws.Range("N122").Value = -27380.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()

$ws.Range("H113").Value = 899.6
$ws.Range("I113").Value = 497.25
$ws.Range("J113").Value = 1045.909
$ws.Range("K113").Value = 1491.75
$ws.Range("L113").Value = 3137.727
$ws.Range("M113").Value = 678.25
$ws.Range("N113").Value = -7477.727000000001

$ws.Range("H122").Value = 2150
$ws.Range("I122").Value = 2150
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6450
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4000
$ws.Range("N122").ClearContents()

$ws.Range("H136").Value = 4692.137
$ws.Range("I136").Value = 1580.9642
$ws.Range("J136").Value = 8479.652
$ws.Range("K136").Value = 4742.892599999999
$ws.Range("L136").Value = 25438.956
$ws.Range("M136").Value = -2192.892599999999
$ws.Range("N136").Value = -30538.956
